# 01.01.2025 - Dodanie sortowania oraz wyrownania danych w tabelach
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "Sortowanie" / wyrownanie danych w kolumnie D - uzupelnienie brakujacych wartosci
$ws.Range("D5").Value = 100
$ws.Range("D6").Value = 100
$ws.Range("D8").Value = 100
$ws.Range("D9").Value = 100

# Przywrocenie widoku arkusza na gore tabeli (reset przewiniecia) i ustawienie
# aktywnej komorki/zaznaczenia tak jak zapisane zostalo przy zamknieciu pliku
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("C11").Select()
